$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("2:6").Insert()
$ws.Range("A2:J6").ClearFormats()
